$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 5 de Mayo de 2020 a las 08:03"

# --- India (row 18) new case numbers ---
$ws.Range("B18").Value = 46476
$ws.Range("C18").Value = 39
$ws.Range("D18").Value = 12849
$ws.Range("E18").Value = 32056
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = 1571

# --- Reorder Bulgaria ahead of Bolivia/Cuba with new Bulgaria figures ---
# Remove the old Bulgaria row (currently row 82, between Cuba and Macedonia)
$ws.Rows("82:82").Delete()
# Insert a new row before Bolivia (row 80) for Bulgaria's updated figures
$ws.Rows("80:80").Insert()
$ws.Range("A80").Value = "Bulgaria"
$ws.Range("B80").Value = 1689
$ws.Range("C80").Value = 37
$ws.Range("D80").Value = 342
$ws.Range("E80").Value = 1269
$ws.Range("F80").Value = 37
$ws.Range("G80").Value = 0
$ws.Range("H80").Value = 78

# --- Maldivas (row 115) updated figures ---
$ws.Range("D115").Value = 17
$ws.Range("E115").Value = 533

# --- Vietnam (row 132) updated figures ---
$ws.Range("D132").Value = 232
$ws.Range("E132").Value = 39
